$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-30, regenerated after
# switching from Strike# to K and recalculating std/mean/s_vals.
$newValues = @{
    2  = 0
    3  = 1
    4  = 4
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 1
    15 = 3
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 4
    29 = 2
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
